# Add a new "Stamina System" bullet to the "Added" list, right after
# the existing "GNU License" bullet (same ListParagraph / numId=3 list).

$d = $word.ActiveDocument

# Locate the "GNU License" paragraph via Find so we don't depend on a
# hard-coded paragraph index.
$range = $d.Content
$found = $range.Find.Execute("GNU License", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

if ($found -and $range.Find.Found) {
    $gnuPara = $range.Paragraphs(1)

    # Insert a new paragraph right after it; Word carries over the
    # paragraph formatting (style + numbering) of the source paragraph.
    $gnuPara.Range.InsertParagraphAfter()

    $newPara = $gnuPara.Next()
    $newPara.Range.Text = "Stamina System"
}
